$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price/Volume columns are treated as text so values like "7.160"
# or "0.00001090" keep their exact literal form instead of being parsed
# as numbers (which would drop trailing/leading zeros).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.013.31"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.42%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.851.32"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -1.10%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9992"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.56%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "311.03"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.82%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9981"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.64%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5068"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -1.33%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3895"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.39%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08225"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -1.70%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.105"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -0.82%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "41.36"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.79%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.183"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.18%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.842.61"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.97%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "20.13"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -2.15%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.160"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -1.83%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.9984"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.60%  "

$ws.Range("B17").Value = "ShibaInu"
$ws.Range("C17").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001090"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -1.41%  "

$ws.Range("B18").Value = "Litecoin"
$ws.Range("C18").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "90.66"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.35%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06635"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.45%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.50"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -1.14%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.9983"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.54%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.898"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -2.23%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "28.039.93"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.44%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.01"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.93%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.225"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -1.13%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.042.93"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -1.83%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "158.78"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.12%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.37"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.98%  "

$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -3.22%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "124.59"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.21%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.1049"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -1.07%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.027"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -1.11%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.791"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -1.68%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.571"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.91%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.02415"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.96%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.06432"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -1.64%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "9.007"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -6.18%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2154"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -1.46%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.240"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +0.76%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.6389"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -1.65%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.173"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -2.63%  "

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -1.71%  "

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -3.02%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.5975"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -1.71%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "12.89"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.84%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.643"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -1.01%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.263"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -1.23%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.990"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.97%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.197"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -1.57%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "120.12"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -1.11%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06842"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.44%  "
